# Fill in results for tasks 2, 3 and 5 (columns C, D, F) on Лист1,
# and clear the old B7 "No file" placeholder (row 7 had no task_1 result).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Header row
$ws.Range("C1").Value = "task_2"
$ws.Range("D1").Value = "task_3"
$ws.Range("F1").Value = "task_5"

# Row 2 - Бикбаев
$ws.Range("C2").Value = 25
$ws.Range("D2").Value = "No file"
$ws.Range("F2").Value = 40

# Row 3 - Вересович
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = "No file"
$ws.Range("F3").Value = 40

# Row 4 - Евланов
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("F4").Value = 40

# Row 5 - Казак
$ws.Range("C5").Value = 25
$ws.Range("D5").Value = 25
$ws.Range("F5").Value = 40

# Row 6 - Капустин
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = "No file"
$ws.Range("F6").Value = 40

# Row 7 - колбов : no task_1 submission, clear the old value and set
# task_2/task_3/task_5 results to "No file"
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = "No file"
$ws.Range("D7").Value = "No file"
$ws.Range("F7").Value = "No file"

# Match the final selection recorded in the workbook - the whole of row 7
$ws.Range("A7:XFD7").Select()
